# Actualización 11 de Mayo - Mañana
# Two more "rescatable" (make-up exam) students are added on the
# "Rescatables" sheet, above the student who was already listed there.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rescatables")

# The student that used to sit in row 2 now moves down to row 4 so the
# two new students can be listed first.
$ws.Cells.Item(4, 1).Value = 19330051920369
$ws.Cells.Item(4, 2).Value = "DE LOS SANTOS"
$ws.Cells.Item(4, 3).Value = "XOTLANIHUA"
$ws.Cells.Item(4, 4).Value = "JENNIFER"
$ws.Cells.Item(4, 5).Value = "ECOLOGÍA"
$ws.Cells.Item(4, 6).Value = "4ARHV"
$ws.Cells.Item(4, 7).Value = 1

# New row 2: COSCAHUA TZOYONTLE, ALMA LIZETH - grupo 4AEV, 2 reprobadas
$ws.Cells.Item(2, 1).Value = 19330051920046
$ws.Cells.Item(2, 2).Value = "COSCAHUA"
$ws.Cells.Item(2, 3).Value = "TZOYONTLE"
$ws.Cells.Item(2, 4).Value = "ALMA LIZETH"
$ws.Cells.Item(2, 5).Value = "ECOLOGÍA"
$ws.Cells.Item(2, 6).Value = "4AEV"
$ws.Cells.Item(2, 7).Value = 2

# New row 3: GONZALEZ SANCHEZ, JONATHAN - grupo 4AEV, 2 reprobadas
$ws.Cells.Item(3, 1).Value = 19330051920057
$ws.Cells.Item(3, 2).Value = "GONZALEZ"
$ws.Cells.Item(3, 3).Value = "SANCHEZ"
$ws.Cells.Item(3, 4).Value = "JONATHAN"
$ws.Cells.Item(3, 5).Value = "ECOLOGÍA"
$ws.Cells.Item(3, 6).Value = "4AEV"
$ws.Cells.Item(3, 7).Value = 2
